# Updated cryptos list (price + 1h volume change) - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking string into a cell while keeping it stored as
# plain text (matches the source data, which is pre-formatted display text,
# not a real number) and without leaving a residual custom style on the cell.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "65.476.99"
$ws.Range("E2").Value = "  -0.75%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.262.74"
$ws.Range("E3").Value = "  -0.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
Set-TextValue "D5" "583.26"
$ws.Range("E5").Value = "  +1.75%  "

# Row 6 - Solana
Set-TextValue "D6" "176.58"
$ws.Range("E6").Value = "  -0.62%  "

# Row 7 - XRP
Set-TextValue "D7" "0.636"
$ws.Range("E7").Value = "  +1.19%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.06%  "

# Row 9 - LidoStakedEther
$ws.Range("D9").Value = "3.266.47"
$ws.Range("E9").Value = "  -0.19%  "

# Row 10 - Dogecoin
Set-TextValue "D10" "0.123"
$ws.Range("E10").Value = "  -1.66%  "

# Row 11 - Toncoin
Set-TextValue "D11" "6.83"
$ws.Range("E11").Value = "  +1.82%  "

# Row 12 - Cardano
Set-TextValue "D12" "0.395"
$ws.Range("E12").Value = "  -0.86%  "

# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.840.61"
$ws.Range("E13").Value = "  -0.12%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -2.66%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "65.669.40"
$ws.Range("E15").Value = "  -0.50%  "

# Row 16 - Avalanche
Set-TextValue "D16" "26.04"
$ws.Range("E16").Value = "  -1.57%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "3.286.26"
$ws.Range("E17").Value = "  +0.49%  "

# Row 18 - ShibaInu
$ws.Range("E18").Value = "  -0.89%  "

# Row 19 - BitcoinCash
Set-TextValue "D19" "419.54"
$ws.Range("E19").Value = "  -3.54%  "

# Row 20 - Polkadot
Set-TextValue "D20" "5.43"
$ws.Range("E20").Value = "  -2.27%  "

# Row 21 - Chainlink
Set-TextValue "D21" "12.89"
$ws.Range("E21").Value = "  -1.91%  "

# Row 22 - Uniswap
Set-TextValue "D22" "7.21"
$ws.Range("E22").Value = "  -2.33%  "

# Row 23 - Dai
$ws.Range("E23").Value = "  -0.11%  "

# Row 24 - Litecoin
Set-TextValue "D24" "70.82"
$ws.Range("E24").Value = "  -1.69%  "

# Row 25 - LEO
Set-TextValue "D25" "5.65"
$ws.Range("E25").Value = "  -0.69%  "

# Row 26 - Kaspa
$ws.Range("E26").Value = "  +5.13%  "

# Row 27 - Polygon
Set-TextValue "D27" "0.501"
$ws.Range("E27").Value = "  -0.65%  "

# Row 28 - PEPE
$ws.Range("E28").Value = "  -0.62%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextValue "D29" "9.31"
$ws.Range("E29").Value = "  +5.14%  "

# Row 31 - PancakeSwap
Set-TextValue "D31" "1.90"
$ws.Range("E31").Value = "  -1.37%  "

# Row 32 - EthereumClassic
Set-TextValue "D32" "22.05"
$ws.Range("E32").Value = "  -1.03%  "

# Row 33 - USDe
$ws.Range("E33").Value = "  -0.02%  "

# Row 34 - NEARProtocol
Set-TextValue "D34" "5.07"
$ws.Range("E34").Value = "  -1.30%  "

# Row 35 - Aptos
Set-TextValue "D35" "6.50"
$ws.Range("E35").Value = "  -1.19%  "

# Row 36 - Fetch.AI
$ws.Range("E36").Value = "  -0.96%  "

# Row 37 - Monero
Set-TextValue "D37" "158.40"
$ws.Range("E37").Value = "  -0.64%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "1.41"
$ws.Range("E38").Value = "  -1.32%  "

# Row 39 - Maker
$ws.Range("D39").Value = "2.833.38"
$ws.Range("E39").Value = "  +2.80%  "

# Row 40 - Stacks
Set-TextValue "D40" "1.76"
$ws.Range("E40").Value = "  -1.58%  "

# Row 41 - EnergySwap
Set-TextValue "D41" "25.80"
$ws.Range("E41").Value = "  -2.98%  "

# Row 42 - Filecoin
Set-TextValue "D42" "4.29"
$ws.Range("E42").Value = "  -0.60%  "

# Row 43 - Mantle
Set-TextValue "D43" "0.742"
$ws.Range("E43").Value = "  -4.28%  "

# Row 44 - OKB
Set-TextValue "D44" "39.48"
$ws.Range("E44").Value = "  -1.80%  "

# Row 45 - RenderToken
Set-TextValue "D45" "5.81"
$ws.Range("E45").Value = "  -3.44%  "

# Row 46 - Hedera
Set-TextValue "D46" "0.0632"
$ws.Range("E46").Value = "  -3.60%  "

# Row 47 - ranking swap: was dogwifhat, now Bittensor
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D47" "310.33"
$ws.Range("E47").Value = "  -3.32%  "

# Row 48 - ranking swap: was Bittensor, now dogwifhat
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D48" "2.22"
$ws.Range("E48").Value = "  -2.97%  "

# Row 49 - InjectiveProtocol
Set-TextValue "D49" "22.46"
$ws.Range("E49").Value = "  -3.59%  "

# Row 50 - VeChain
Set-TextValue "D50" "0.0265"
$ws.Range("E50").Value = "  -0.49%  "

# Row 51 - Stellar
Set-TextValue "D51" "0.101"
$ws.Range("E51").Value = "  -0.88%  "

